# Update for release to deploy 0.1.1
$wb = $excel.ActiveWorkbook

# 1. Rename the "Include from LOINC" sheet to "Include #0"
$wsInclude = $wb.Worksheets.Item("Include from LOINC")
$wsInclude.Name = "Include #0"

# 2. Update the Metadata sheet
$wsMeta = $wb.Worksheets.Item("Metadata")

# Version 0.1.0 -> 0.1.1
$wsMeta.Range("B3").Value = "0.1.1"

# Date bump
$wsMeta.Range("B8").Value = "2024-11-11T17:53:38-06:00"

# 3. Insert a new "Jurisdiction" property row right before "Description" (row 11),
#    pushing Description/Purpose/Copyright/Immutable down by one row.
$wsMeta.Rows.Item(11).Insert()

# Copy formatting from the row that is now right below (old row 11, now row 12)
# so the new row matches the existing table styling.
$wsMeta.Range("A12:B12").Copy($wsMeta.Range("A11:B11"))

# Set the new row's actual content.
$wsMeta.Range("A11").Value = "Jurisdiction"
$wsMeta.Range("B11").Value = ""
